# no-op test
$p = $ppt.ActivePresentation
$n = $p.Slides.Count
Write-Host "Slide count: $n"
